$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 189, shifting existing rows 189:219 down to 190:220
$ws.Rows(189).EntireRow.Insert()

# Populate the new row 189 with the new record
$ws.Cells.Item(189, 1).Value = 10
$ws.Cells.Item(189, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(189, 3).Value = "La Araucanía"
$ws.Cells.Item(189, 4).Value = 44504
$ws.Cells.Item(189, 5).Value = 9
$ws.Cells.Item(189, 6).Value = 100114013
$ws.Cells.Item(189, 7).Value = "Zanahoria"
$ws.Cells.Item(189, 8).Value = "Sin especificar"
$ws.Cells.Item(189, 9).Value = "Primera"
$ws.Cells.Item(189, 10).Value = 210
$ws.Cells.Item(189, 11).Value = 7000
$ws.Cells.Item(189, 12).Value = 7000
$ws.Cells.Item(189, 13).Value = 7000
$ws.Cells.Item(189, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(189, 15).Value = "Región del Maule"
$ws.Cells.Item(189, 16).Value = 350
$ws.Cells.Item(189, 17).Value = 20
$ws.Cells.Item(189, 18).Value = "Hortaliza"
